# Updated cryptos list on Mon Feb 19 16:13:15 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.139.07"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").Value = "2.910.96"
$ws.Range("E3").Value = "  +3.75%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "355.26"
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "113.86"
$ws.Range("E6").Value = "  +1.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.559"
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.623"
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.66"
$ws.Range("E10").Value = "  -1.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0877"
$ws.Range("E11").Value = "  +4.53%  "
$ws.Range("E12").Value = "  +0.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.80"
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.74"
$ws.Range("E14").Value = "  -0.36%  "
$ws.Range("D15").Value = "3.369.33"
$ws.Range("E15").Value = "  +3.87%  "
$ws.Range("D16").Value = "2.903.83"
$ws.Range("E16").Value = "  +3.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.987"
$ws.Range("E17").Value = "  +3.18%  "
$ws.Range("D18").Value = "52.210.07"
$ws.Range("E18").Value = "  +0.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.34"
$ws.Range("E19").Value = "  +0.62%  "
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.09"
$ws.Range("E21").Value = "  +4.15%  "
$ws.Range("D22").Value = "0.0₃0983"
$ws.Range("E22").Value = "  +0.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.19"
$ws.Range("E23").Value = "  +1.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.62"
$ws.Range("E24").Value = "  +0.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.81"
$ws.Range("E25").Value = "  +2.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.182"
$ws.Range("E26").Value = "  +12.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.81"
$ws.Range("E27").Value = "  +2.47%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.68"
$ws.Range("E29").Value = "  +2.50%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.85"
$ws.Range("E30").Value = "  +12.19%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.103"
$ws.Range("E31").Value = "  +14.07%  "
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "37.53"
$ws.Range("E32").Value = "  -4.07%  "
$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.27"
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.05"
$ws.Range("E34").Value = "  +8.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "53.16"
$ws.Range("E35").Value = "  +1.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0453"
$ws.Range("E36").Value = "  +0.53%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.35"
$ws.Range("E38").Value = "  +5.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.81"
$ws.Range("E39").Value = "  -1.04%  "
$ws.Range("E40").Value = "  +1.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.71"
$ws.Range("E41").Value = "  +7.87%  "
$ws.Range("E42").Value = "  +0.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.00"
$ws.Range("E43").Value = "  +4.65%  "
$ws.Range("E44").Value = "  -2.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "117.88"
$ws.Range("E45").Value = "  -2.21%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.53"
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.51"
$ws.Range("E47").Value = "  +1.28%  "
$ws.Range("D48").Value = "2.179.49"
$ws.Range("E48").Value = "  +2.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.254"
$ws.Range("E49").Value = "  +15.50%  "
$ws.Range("E50").Value = "  +12.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.952"
$ws.Range("E51").Value = "  -2.69%  "
